# Apply cryptos list update (price/volume refresh + BitcoinCash/Polkadot row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.078.46'
$ws.Range("E2").Value = '  +3.23%  '
$ws.Range("D3").Value = '2.421.93'
$ws.Range("E3").Value = '  +4.03%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.83'
$ws.Range("E5").Value = '  +2.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.56'
$ws.Range("E6").Value = '  +5.75%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.532'
$ws.Range("E8").Value = '  +2.54%  '
$ws.Range("D9").Value = '2.420.82'
$ws.Range("E9").Value = '  +4.23%  '
$ws.Range("E10").Value = '  +5.08%  '
$ws.Range("E11").Value = '  +1.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.39'
$ws.Range("E12").Value = '  +2.30%  '
$ws.Range("E13").Value = '  +4.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.33'
$ws.Range("E14").Value = '  +8.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000175'
$ws.Range("E15").Value = '  +10.38%  '
$ws.Range("D16").Value = '2.845.82'
$ws.Range("E16").Value = '  +3.64%  '
$ws.Range("D17").Value = '61.737.34'
$ws.Range("E17").Value = '  +2.56%  '
$ws.Range("D18").Value = '2.418.48'
$ws.Range("E18").Value = '  +3.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.12'
$ws.Range("E19").Value = '  +6.22%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '324.77'
$ws.Range("E20").Value = '  +3.36%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.19'
$ws.Range("E21").Value = '  +3.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.77'
$ws.Range("E22").Value = '  +4.47%  '
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.91'
$ws.Range("E24").Value = '  +3.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.76'
$ws.Range("E25").Value = '  +5.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.17'
$ws.Range("E26").Value = '  +11.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '561.61'
$ws.Range("E27").Value = '  +14.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.60%  '
$ws.Range("D29").Value = '2.513.55'
$ws.Range("E29").Value = '  +3.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.36'
$ws.Range("E30").Value = '  +6.70%  '
$ws.Range("E31").Value = '  +9.43%  '
$ws.Range("E32").Value = '  +6.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.148'
$ws.Range("E33").Value = '  +3.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.86'
$ws.Range("E34").Value = '  +4.76%  '
$ws.Range("E35").Value = '  +4.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.77'
$ws.Range("E36").Value = '  +12.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.98'
$ws.Range("E37").Value = '  +12.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.87%  '
$ws.Range("E39").Value = '  +6.72%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.384'
$ws.Range("E40").Value = '  +3.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.78'
$ws.Range("E41").Value = '  +2.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '146.77'
$ws.Range("E42").Value = '  +3.94%  '
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.30'
$ws.Range("E44").Value = '  +13.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '149.52'
$ws.Range("E45").Value = '  +6.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.63'
$ws.Range("E46").Value = '  +3.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0539'
$ws.Range("E47").Value = '  +6.41%  '
$ws.Range("E48").Value = '  +7.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.590'
$ws.Range("E49").Value = '  +4.75%  '
$ws.Range("E50").Value = '  +4.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0909'
$ws.Range("E51").Value = '  +2.15%  '
